$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.991708755493164
$ws.Range("B1").Value = 2.785296678543091
$ws.Range("C1").Value = 1.600721120834351
$ws.Range("D1").Value = 1.250567197799683
$ws.Range("E1").Value = 1.135878324508667
